# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E) and "Correspond Handback
# DateTime" (H) columns for the e460d230-... row (row 3) on the "zh-cn"
# sheet, and for the e460d230-... row (row 3) on the "de-de" sheet, to
# reflect the new handback timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-20 02:36:44"
$wsZhCn.Range("H3").Value = "2016-03-20 02:37:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-20 02:36:47"
$wsDeDe.Range("H3").Value = "2016-03-20 02:37:08"
